# edit.ps1 -- applies the "finished draft hydro metrics table in methods" commit
#
# Summary of changes:
#  1. "Thus plot mean wood density is maximised when flows patterns..." gains
#     the word "average" and two new sentences about wood density being
#     negatively predicted by interannual uniformity of minimum flows.
#  2. The following paragraph ("A similar relationship...") loses the word
#     "significantly" and gains a "base" prefix turning "flow deviated" into
#     "baseflow deviated"; the document's "_GoBack" bookmark moves here
#     (between "base" and "flow").
#  3. "Metrics of low flow duration were not significantly predictive of
#     wood density." is reworded to "...did not significantly predict wood
#     density."
#  4. The old "_GoBack" bookmark (previously sitting before "dense stemmed")
#     is removed -- handled implicitly because bookmark names are unique:
#     adding "_GoBack" at its new location removes it from the old one.

$d = $word.ActiveDocument

# --- 1a. "maximised when flows patterns" -> "maximised when average flows patterns"
$d.Content.Find.Execute(
    "maximised when flows patterns", $true, $false, $false, $false, $false,
    $true, 1, $false, "maximised when average flows patterns", 2) | Out-Null

# --- 1b. Insert the two new sentences at the end of that same paragraph,
#         right after "...is not consistent throughout the record. "
$rng = $d.Content
$rng.Find.Execute(
    "is not consistent throughout the record. ", $true, $false, $false,
    $false, $false, $true, 1, $false, $null, 0) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)
$newSentences = "Wood density was negatively predicted by interannual uniformity (contingency), but not constancy of minimum flows. That is to say, it was not important how strongly minimum flows were associated with particular seasons, but whether the seasonal pattern of flows was the same across years of the record.  "
$insertPoint.InsertAfter($newSentences)

# --- 2a. Drop "significantly " from "flow deviated significantly from the mean."
$d.Content.Find.Execute(
    "flow deviated significantly from the mean", $true, $false, $false,
    $false, $false, $true, 1, $false, "flow deviated from the mean", 2) | Out-Null

# --- 2b. Turn "...years in which flow deviated..." into
#         "...years in which baseflow deviated..." and drop the "_GoBack"
#         bookmark right in the middle of the new word ("base" | "flow").
$rng2 = $d.Content
$rng2.Find.Execute(
    "years in which ", $true, $false, $false, $false, $false, $true, 1,
    $false, $null, 0) | Out-Null
$baseInsertPoint = $d.Range($rng2.End, $rng2.End)
$baseInsertPoint.InsertAfter("base")
$bookmarkRange = $d.Range($rng2.End + 4, $rng2.End + 4)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# --- 3. Reword the low-flow-duration sentence.
$d.Content.Find.Execute(
    "Metrics of low flow duration were not significantly predictive of wood density.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Metrics of low flow duration did not significantly predict wood density.", 2) | Out-Null
